$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: append a new row 14 that is a duplicate of the current row 13
# (the weekly entry that is about to be superseded by fresher data).
# Value2 is used because it reliably round-trips the underlying cell
# value/type through this COM layer (unlike Value in this runtime).
$ws.Range("A14").Value2 = $ws.Range("A13").Value2
$ws.Range("B14").Value2 = $ws.Range("B13").Value2
$ws.Range("C14").Value2 = $ws.Range("C13").Value2
$ws.Range("D14").Value2 = $ws.Range("D13").Value2
$ws.Range("D14").NumberFormat = $ws.Range("D13").NumberFormat
$ws.Range("E14").Value2 = $ws.Range("E13").Value2
$ws.Range("F14").Value2 = $ws.Range("F13").Value2
$ws.Range("G14").Value2 = $ws.Range("G13").Value2
$ws.Range("H14").Value2 = $ws.Range("H13").Value2
$ws.Range("I14").Value2 = $ws.Range("I13").Value2
$ws.Range("J14").Value2 = $ws.Range("J13").Value2
$ws.Range("K14").Value2 = $ws.Range("K13").Value2
$ws.Range("L14").Value2 = $ws.Range("L13").Value2
$ws.Range("M14").Value2 = $ws.Range("M13").Value2
$ws.Range("N14").Value2 = $ws.Range("N13").Value2
$ws.Range("O14").Value2 = $ws.Range("O13").Value2
$ws.Range("P14").Value2 = $ws.Range("P13").Value2
$ws.Range("Q14").Value2 = $ws.Range("Q13").Value2
$ws.Range("R14").Value2 = $ws.Range("R13").Value2

# Step 2: update row 13 in place with the newer weekly price data.
$ws.Range("D13").Value2 = 45142
$ws.Range("K13").Value2 = 17000
$ws.Range("L13").Value2 = 18000
$ws.Range("M13").Value2 = 17500
$ws.Range("P13").Value2 = 972
